$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6727.6665
$ws.Range("I62").Value = 5107
$ws.Range("J62").Value = 9969
$ws.Range("K62").Value = 5107
$ws.Range("L62").Value = 9969
$ws.Range("M62").Value = -4483
$ws.Range("N62").Value = -11217

$ws.Range("H65").Value = 6727.6665
$ws.Range("I65").Value = 5107
$ws.Range("J65").Value = 9969
$ws.Range("K65").Value = 25535
$ws.Range("L65").Value = 49845
$ws.Range("M65").Value = -22415
$ws.Range("N65").Value = -56085

$ws.Range("H76").Value = 7749.75
$ws.Range("J76").Value = 8666.333000000001
$ws.Range("L76").Value = 8666.333000000001
$ws.Range("N76").Value = -9296.333000000001

$ws.Range("H79").Value = 7749.75
$ws.Range("J79").Value = 8666.333000000001
$ws.Range("L79").Value = 8666.333000000001
$ws.Range("N79").Value = -10850.333

$ws.Range("H98").Value = 4179.6665
$ws.Range("I98").Value = 4179.6665
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4179.6665
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -2681.6665

$ws.Range("H112").Value = 30753.143
$ws.Range("J112").Value = 33514.688
$ws.Range("L112").Value = 100544.064
$ws.Range("N112").Value = -102760.064

$ws.Range("H122").Value = 4179.6665
$ws.Range("I122").Value = 4179.6665
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12538.9995
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10088.9995
$ws.Range("N122").ClearContents()

$ws.Range("H125").Value = 3949.5
$ws.Range("I125").Value = 3949.5
$ws.Range("K125").Value = 35545.5
$ws.Range("M125").Value = -33085.5

$ws.Range("H138").Value = 6026734.5
$ws.Range("I138").Value = 1402.25
$ws.Range("J138").Value = 8477717
$ws.Range("K138").Value = 4206.75
$ws.Range("L138").Value = 25433151
$ws.Range("M138").Value = 933.25
$ws.Range("N138").Value = -25443431

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 11000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H74").Value = 52381.227
$ws.Range("I74").Value = 68752.87
$ws.Range("J74").Value = 17299.143
$ws.Range("K74").Value = 68752.87
$ws.Range("L74").Value = 17299.143
$ws.Range("M74").Value = -67878.87
$ws.Range("N74").Value = -19047.143

$ws.Range("H77").Value = 52381.227
$ws.Range("I77").Value = 68752.87
$ws.Range("J77").Value = 17299.143
$ws.Range("K77").Value = 343764.35
$ws.Range("L77").Value = 86495.715
$ws.Range("M77").Value = -339396.35
$ws.Range("N77").Value = -95231.715

$ws.Range("H122").Value = 2105.3635
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2332.375
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 6997.125
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -11897.125

$ws.Range("H125").Value = 24392
$ws.Range("J125").Value = 24392
$ws.Range("L125").Value = 24392
$ws.Range("N125").Value = -34232

$ws.Range("H132").Value = 2845.5144
$ws.Range("I132").Value = 2456.375
$ws.Range("K132").Value = 7369.125
$ws.Range("M132").Value = -4839.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 15920.143
$ws.Range("I26").Value = 15920.143
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 15920.143
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -15628.143

$ws.Range("H134").Value = 1895.0883
$ws.Range("I134").Value = 1893.5223
$ws.Range("K134").Value = 5680.5669
$ws.Range("M134").Value = -3145.5669

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50837.24
$ws.Range("I31").Value = 69546.664
$ws.Range("J31").Value = 4063.6667
$ws.Range("K31").Value = 69546.664
$ws.Range("L31").Value = 4063.6667
$ws.Range("M31").Value = -69251.664
$ws.Range("N31").Value = -4653.6667

$ws.Range("H34").Value = 50837.24
$ws.Range("I34").Value = 69546.664
$ws.Range("J34").Value = 4063.6667
$ws.Range("K34").Value = 69546.664
$ws.Range("L34").Value = 4063.6667
$ws.Range("M34").Value = -69344.664
$ws.Range("N34").Value = -4467.6667

$ws.Range("H60").Value = 16634
$ws.Range("J60").Value = 14951.5
$ws.Range("L60").Value = 14951.5
$ws.Range("N60").Value = -15973.5

$ws.Range("H99").Value = 3605.4443
$ws.Range("J99").Value = 3838.25
$ws.Range("L99").Value = 3838.25
$ws.Range("N99").Value = -6834.25

$ws.Range("H126").Value = 3605.4443
$ws.Range("J126").Value = 3838.25
$ws.Range("L126").Value = 11514.75
$ws.Range("N126").Value = -16454.75

$ws.Range("H132").Value = 4079.7856
$ws.Range("I132").Value = 4009
$ws.Range("K132").Value = 12027
$ws.Range("M132").Value = -9497

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1071.1111
$ws.Range("I26").Value = 380.81818
$ws.Range("J26").Value = 2155.8572
$ws.Range("K26").Value = 1142.45454
$ws.Range("L26").Value = 6467.571599999999
$ws.Range("M26").Value = -854.45454
$ws.Range("N26").Value = -7043.571599999999

$ws.Range("H76").Value = 2250
$ws.Range("I76").Value = 2250
$ws.Range("K76").Value = 6750
$ws.Range("M76").Value = -6367

$ws.Range("H79").Value = 2250
$ws.Range("I79").Value = 2250
$ws.Range("K79").Value = 6750
$ws.Range("M79").Value = -5424

$ws.Range("H131").Value = 31623.53
$ws.Range("J131").Value = 2037.5
$ws.Range("L131").Value = 6112.5
$ws.Range("N131").Value = -16192.5

$ws.Range("H132").Value = 1540.5264
$ws.Range("J132").Value = 1929.6666
$ws.Range("L132").Value = 17366.9994
$ws.Range("N132").Value = -22426.9994

$ws.Range("H140").Value = 3465.6667
$ws.Range("I140").Value = 1397
$ws.Range("J140").Value = 4500
$ws.Range("K140").Value = 4191
$ws.Range("L140").Value = 13500
$ws.Range("M140").Value = 989
$ws.Range("N140").Value = -23860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 13000
$ws.Range("J21").Value = 18500
$ws.Range("L21").Value = 18500
$ws.Range("N21").Value = -18846

$ws.Range("H30").Value = 13000
$ws.Range("J30").Value = 18500
$ws.Range("L30").Value = 18500
$ws.Range("N30").Value = -18710

$ws.Range("H102").Value = 43480376
$ws.Range("I102").Value = 1928.1052
$ws.Range("K102").Value = 1928.1052
$ws.Range("M102").Value = -306.1052

$ws.Range("H122").Value = 3247.5
$ws.Range("I122").Value = 2853.647
$ws.Range("K122").Value = 8560.940999999999
$ws.Range("M122").Value = -6110.940999999999

$ws.Range("H130").Value = 58792.668
$ws.Range("J130").Value = 58792.668
$ws.Range("L130").Value = 58792.668
$ws.Range("N130").Value = -68832.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 99985.5
$ws.Range("J121").Value = 99985.5
$ws.Range("L121").Value = 99985.5
$ws.Range("N121").Value = -103479.5

$ws.Range("H132").Value = 4165.7144
$ws.Range("I132").Value = 3892.9285
$ws.Range("J132").Value = 4711.2856
$ws.Range("K132").Value = 11678.7855
$ws.Range("L132").Value = 14133.8568
$ws.Range("M132").Value = -9148.7855
$ws.Range("N132").Value = -19193.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 8935
$ws.Range("I51").Value = 8935
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 8935
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -8425

$ws.Range("H70").Value = 10095
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 10095
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H122").Value = 2399.96
$ws.Range("I122").Value = 2045.4375
$ws.Range("J122").Value = 3030.2222
$ws.Range("K122").Value = 6136.3125
$ws.Range("L122").Value = 9090.6666
$ws.Range("M122").Value = -3686.3125
$ws.Range("N122").Value = -13990.6666
